$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old F23 content ("test hp ...") moves down to the new F24 cell, carrying
# its original formatting (style s="5": plain Calibri, m/d/yyyy, left aligned).
$ws.Range("F23").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("F24").Value = $ws.Range("F23").Value2

# F23 gets the new bold "recap" text (reworded hyp-test row, beta-blocker/ROPE example).
$ws.Range("F23").Value = "recap test 2 campioni; it- test bayesiano"
$ws.Range("F23").Font.Bold = $true

# E24 continues the E15:E23 "+7 days" weekly date series.
$ws.Range("E24").Formula = "=E22+7"
